$d = $word.ActiveDocument

# The authored change renames the inline picture shapes that live in the
# document's headers/footers:
#   - the Pearson logo (PNG) is renamed to "image1.png"
#   - the BTEC logo (JPG) is renamed to "image2.jpg"
# `InlineShape.Name` is write-only in this host (the getter always reads
# back empty), so identify which picture we are looking at from its
# AlternativeText (the descr attribute), which *is* readable.

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $n = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $n; $j++) {
                $ishp = $hdr.Range.InlineShapes.Item($j)
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $ishp.Name = "image2.jpg"
                } elseif ($ishp.AlternativeText -like "*PearsonLogo.png") {
                    $ishp.Name = "image1.png"
                }
            }
        }
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $n = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $n; $j++) {
                $ishp = $ftr.Range.InlineShapes.Item($j)
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $ishp.Name = "image2.jpg"
                } elseif ($ishp.AlternativeText -like "*PearsonLogo.png") {
                    $ishp.Name = "image1.png"
                }
            }
        }
    }
}

Write-Output "renamed inline shapes"
